$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.026972333333333
$ws.Range("H2").Value = 3.080917
$ws.Range("I2").Value = 0.2032541865322035
$ws.Range("J2").Value = 0.2032541865322035
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 42.23928166666666
$ws.Range("N2").Value = 126.717845
$ws.Range("O2").Value = 0.7297675404946526
$ws.Range("P2").Value = 0.7297675404946528
$ws.Range("Q2").Value = 43.37857365154056
$ws.Range("R2").Value = 390.407162863865
$ws.Range("S2").Value = 0.1483283078008475
$ws.Range("T2").Value = 0.1483283078008475

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.026972333333333
$ws.Range("H3").Value = 3.080917
$ws.Range("I3").Value = 0.2032541865322035
$ws.Range("J3").Value = 0.2032541865322035
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.624984
$ws.Range("N3").Value = 13.874952
$ws.Range("O3").Value = 0.07990579066051323
$ws.Range("P3").Value = 0.07990579066051323
$ws.Range("Q3").Value = 4.749730610109334
$ws.Range("R3").Value = 42.747575490984
$ws.Range("S3").Value = 0.01624118647991516
$ws.Range("T3").Value = 0.01624118647991516

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.026972333333333
$ws.Range("H4").Value = 3.080917
$ws.Range("I4").Value = 0.2032541865322035
$ws.Range("J4").Value = 0.2032541865322035
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.151907
$ws.Range("N4").Value = 6.455721
$ws.Range("O4").Value = 0.03717847029587412
$ws.Range("P4").Value = 0.03717847029587412
$ws.Range("Q4").Value = 2.209948952906334
$ws.Range("R4").Value = 19.889540576157
$ws.Range("S4").Value = 0.007556679736499584
$ws.Range("T4").Value = 0.007556679736499584

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.026972333333333
$ws.Range("H5").Value = 3.080917
$ws.Range("I5").Value = 0.2032541865322035
$ws.Range("J5").Value = 0.2032541865322035
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.864288333333333
$ws.Range("N5").Value = 26.592865
$ws.Range("O5").Value = 0.1531481985489599
$ws.Range("P5").Value = 0.15314819854896
$ws.Range("Q5").Value = 9.103378873022777
$ws.Range("R5").Value = 81.93040985720499
$ws.Range("S5").Value = 0.03112801251494124
$ws.Range("T5").Value = 0.03112801251494125

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.358031333333334
$ws.Range("H6").Value = 4.074094000000001
$ws.Range("I6").Value = 0.2687760370778347
$ws.Range("J6").Value = 0.2687760370778347
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 42.23928166666666
$ws.Range("N6").Value = 126.717845
$ws.Range("O6").Value = 0.7297675404946526
$ws.Range("P6").Value = 0.7297675404946528
$ws.Range("Q6").Value = 57.36226800082556
$ws.Range("R6").Value = 516.2604120074301
$ws.Range("S6").Value = 0.196144027522191
$ws.Range("T6").Value = 0.196144027522191

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.358031333333334
$ws.Range("H7").Value = 4.074094000000001
$ws.Range("I7").Value = 0.2687760370778347
$ws.Range("J7").Value = 0.2687760370778347
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.624984
$ws.Range("N7").Value = 13.874952
$ws.Range("O7").Value = 0.07990579066051323
$ws.Range("P7").Value = 0.07990579066051323
$ws.Range("Q7").Value = 6.280873188165335
$ws.Range("R7").Value = 56.52785869348801
$ws.Range("S7").Value = 0.0214767617533038
$ws.Range("T7").Value = 0.0214767617533038

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.358031333333334
$ws.Range("H8").Value = 4.074094000000001
$ws.Range("I8").Value = 0.2687760370778347
$ws.Range("J8").Value = 0.2687760370778347
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.151907
$ws.Range("N8").Value = 6.455721
$ws.Range("O8").Value = 0.03717847029587412
$ws.Range("P8").Value = 0.03717847029587412
$ws.Range("Q8").Value = 2.922357132419334
$ws.Range("R8").Value = 26.301214191774
$ws.Range("S8").Value = 0.009992681910741036
$ws.Range("T8").Value = 0.009992681910741036

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.358031333333334
$ws.Range("H9").Value = 4.074094000000001
$ws.Range("I9").Value = 0.2687760370778347
$ws.Range("J9").Value = 0.2687760370778347
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.864288333333333
$ws.Range("N9").Value = 26.592865
$ws.Range("O9").Value = 0.1531481985489599
$ws.Range("P9").Value = 0.15314819854896
$ws.Range("Q9").Value = 12.03798130436778
$ws.Range("R9").Value = 108.34183173931
$ws.Range("S9").Value = 0.04116256589159884
$ws.Range("T9").Value = 0.04116256589159885

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.667646666666666
$ws.Range("H10").Value = 8.002939999999999
$ws.Range("I10").Value = 0.5279697763899619
$ws.Range("J10").Value = 0.5279697763899619
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 42.23928166666666
$ws.Range("N10").Value = 126.717845
$ws.Range("O10").Value = 0.7297675404946526
$ws.Range("P10").Value = 0.7297675404946528
$ws.Range("Q10").Value = 112.6794789404778
$ws.Range("R10").Value = 1014.1153104643
$ws.Range("S10").Value = 0.3852952051716142
$ws.Range("T10").Value = 0.3852952051716143

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.667646666666666
$ws.Range("H11").Value = 8.002939999999999
$ws.Range("I11").Value = 0.5279697763899619
$ws.Range("J11").Value = 0.5279697763899619
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.624984
$ws.Range("N11").Value = 13.874952
$ws.Range("O11").Value = 0.07990579066051323
$ws.Range("P11").Value = 0.07990579066051323
$ws.Range("Q11").Value = 12.33782315098667
$ws.Range("R11").Value = 111.04040835888
$ws.Range("S11").Value = 0.04218784242729428
$ws.Range("T11").Value = 0.04218784242729428

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.667646666666666
$ws.Range("H12").Value = 8.002939999999999
$ws.Range("I12").Value = 0.5279697763899619
$ws.Range("J12").Value = 0.5279697763899619
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.151907
$ws.Range("N12").Value = 6.455721
$ws.Range("O12").Value = 0.03717847029587412
$ws.Range("P12").Value = 0.03717847029587412
$ws.Range("Q12").Value = 5.740527535526666
$ws.Range("R12").Value = 51.66474781973999
$ws.Range("S12").Value = 0.0196291086486335
$ws.Range("T12").Value = 0.0196291086486335

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.667646666666666
$ws.Range("H13").Value = 8.002939999999999
$ws.Range("I13").Value = 0.5279697763899619
$ws.Range("J13").Value = 0.5279697763899619
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 8.864288333333333
$ws.Range("N13").Value = 26.592865
$ws.Range("O13").Value = 0.1531481985489599
$ws.Range("P13").Value = 0.15314819854896
$ws.Range("Q13").Value = 23.64678922478888
$ws.Range("R13").Value = 212.8211030231
$ws.Range("S13").Value = 0.08085762014241987
$ws.Range("T13").Value = 0.08085762014241989

